$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F19").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("E21").Value = 5
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 0

$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4122)

$ws.Range("A25").Value = 43915
$ws.Range("B25").Value = 900
$ws.Range("C25").Value = 133
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = 115
$ws.Range("F25").Value = 0

$ws.Range("E26").Select()
